$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter..." text and the
# paragraph right after it (the copyright notice), plus the blank paragraph
# that immediately precedes them, then remove that whole block. This leaves
# the "LOM3071: Tratamento de Minerios (Requisito fraco)" paragraph followed
# directly by the existing trailing blank paragraph / page-break paragraph.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text

    if ($text -like "*Ver no Jupiter*") {
        # the blank paragraph immediately before this one is the start of
        # the block to remove
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($text -like "*Powered by Jekyll*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
